# Version 09-25 - Large changes to qoq and ifoCAST evaluation
# Updates the "latest release" absolute GDP series (column B, rows 2-148)
# with revised values, and appends a new quarterly observation in row 149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 61.68613450127155
    3 = 62.42644814919469
    4 = 63.06671941226335
    5 = 63.326829612885
    6 = 63.89707120655552
    7 = 64.64738909296412
    8 = 65.98795705001413
    9 = 66.30809268154846
    10 = 67.70868606951117
    11 = 68.90919468776491
    12 = 70.81
    14 = 70.36
    15 = 71.37
    16 = 72.34999999999999
    18 = 71.73
    19 = 71.53
    20 = 70.95999999999999
    22 = 71.39
    23 = 71.33
    24 = 72.34
    25 = 72.73999999999999
    26 = 73.17
    27 = 73.98
    28 = 73.69
    29 = 74.28
    30 = 74.43000000000001
    31 = 74.45999999999999
    32 = 73.93000000000001
    34 = 75.23
    35 = 75.92
    36 = 75.53
    38 = 76.67
    39 = 77.23999999999999
    40 = 77.98
    41 = 77.62
    42 = 78
    43 = 77.98
    44 = 78.88
    45 = 78.78
    46 = 79.86
    47 = 80.23999999999999
    49 = 82.13
    50 = 82.19
    51 = 81.84
    52 = 83.58
    53 = 83.37
    54 = 83.37
    55 = 83.11
    57 = 83.04000000000001
    58 = 83.61
    59 = 83.38
    60 = 82.25
    61 = 82.37
    62 = 83.12
    63 = 83.22
    65 = 83.58
    66 = 83.31
    67 = 83.28
    68 = 83.41
    69 = 83.91
    70 = 84.56999999999999
    71 = 84.89
    72 = 85.84
    73 = 87.31
    74 = 87.97
    75 = 89.26000000000001
    76 = 89.38
    77 = 89.97
    78 = 90.43000000000001
    79 = 91.09999999999999
    80 = 91.63
    81 = 91.3
    82 = 90.8
    84 = 85.2
    85 = 85.41
    86 = 85.94
    88 = 87.23
    95 = 92.95
    96 = 93.15000000000001
    98 = 93.44
    99 = 93.16
    100 = 92.64
    101 = 93.72
    102 = 94.2
    103 = 94.34
    104 = 95.33
    105 = 95.31
    106 = 95.84
    107 = 96.56999999999999
    108 = 96.36
    109 = 96.86
    110 = 97.38
    111 = 97.88
    112 = 98.73999999999999
    113 = 98.95999999999999
    114 = 99.31
    115 = 99.77
    116 = 101.06
    117 = 101.78
    118 = 102.59
    119 = 103.63
    120 = 103.06
    121 = 103.87
    122 = 103.19
    123 = 103.65
    124 = 104.33
    125 = 104.39
    126 = 104.79
    127 = 104.44
    128 = 102.32
    129 = 93.23999999999999
    130 = 101.33
    131 = 102.31
    132 = 101.68
    133 = 104.07
    134 = 104.16
    135 = 104.72
    136 = 105.43
    137 = 105.59
    138 = 105.9
    139 = 105.53
    140 = 105.03
    141 = 104.95
    142 = 104.95
    143 = 104.66
    144 = 104.55
    145 = 104.28
    146 = 104.3
    147 = 104.49
    148 = 104.81
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $values[$row]
}

# Add new row 149
$ws.Cells.Item(149, 1).Value = 45748
$ws.Cells.Item(149, 2).Value = 104.52

# Copy date style/format from A148 to A149
$ws.Cells.Item(148, 1).Copy() | Out-Null
$ws.Cells.Item(149, 1).PasteSpecial(-4122) | Out-Null

Write-Output "done"